# Update the cryptocurrency price/volume table on Sheet1.
# Column layout: A=index, B=Coin, C=Link, D=Price, E=Volume(1h)
#
# Price values are stored as plain text in the workbook (e.g. "3.309.10" is not
# a real number). When a new price string would otherwise be auto-recognized
# by Excel as a numeric literal (e.g. "568.28"), the cell's number format is
# forced to Text ("@") first so the literal text is preserved verbatim instead
# of being converted into a floating point number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceText {
    param($row, $value)
    # Force the cell to Text format so number-like strings aren't reinterpreted.
    $ws.Cells.Item($row, 4).NumberFormat = "@"
    $ws.Cells.Item($row, 4).Value = $value
}

function Set-Price {
    param($row, $value)
    $ws.Cells.Item($row, 4).Value = $value
}

function Set-Volume {
    param($row, $value)
    $ws.Cells.Item($row, 5).Value = "  $value  "
}

# Row 2 - Bitcoin
Set-Price 2 "60.963.73"
Set-Volume 2 "-5.18%"

# Row 3 - Ethereum
Set-Price 3 "3.311.26"
Set-Volume 3 "-5.27%"

# Row 4 - TetherUSD
Set-Volume 4 "-0.05%"

# Row 5 - BNB
Set-PriceText 5 "568.28"
Set-Volume 5 "-3.55%"

# Row 6 - Solana
Set-PriceText 6 "126.41"
Set-Volume 6 "-5.50%"

# Row 7 - USDC
Set-Volume 7 "-0.06%"

# Row 8 - LidoStakedEther
Set-Price 8 "3.310.48"
Set-Volume 8 "-5.27%"

# Row 9 - XRP
Set-Volume 9 "-2.18%"

# Row 10 - Toncoin
Set-PriceText 10 "7.27"
Set-Volume 10 "-4.68%"

# Row 11 - Dogecoin
Set-Volume 11 "-6.16%"

# Row 12 - Cardano
Set-Volume 12 "-4.17%"

# Row 13 - WrappedliquidstakedEther2.0
Set-Price 13 "3.873.18"
Set-Volume 13 "-5.32%"

# Row 14 - TRON
Set-Volume 14 "-1.44%"

# Row 15 - WrappedEther
Set-Price 15 "3.310.51"
Set-Volume 15 "-5.30%"

# Row 16 - ShibaInu
Set-Volume 16 "-7.48%"

# Row 17 - WrappedBTC
Set-Price 17 "61.025.17"
Set-Volume 17 "-5.02%"

# Row 18 - Avalanche
Set-PriceText 18 "24.26"
Set-Volume 18 "-4.16%"

# Row 19 - Polkadot
Set-Volume 19 "-3.76%"

# Row 20 - Uniswap
Set-Volume 20 "-10.23%"

# Row 21 - Chainlink
Set-PriceText 21 "13.10"
Set-Volume 21 "-3.25%"

# Row 22 - BitcoinCash
Set-PriceText 22 "350.13"
Set-Volume 22 "-9.42%"

# Row 23 - Polygon
Set-PriceText 23 "0.551"
Set-Volume 23 "-5.17%"

# Row 24 - Dai
Set-PriceText 24 "1.00"
Set-Volume 24 "+0.03%"

# Row 25 - WrappedeETH
Set-Price 25 "3.441.87"
Set-Volume 25 "-5.31%"

# Row 26 - Litecoin
Set-PriceText 26 "69.74"
Set-Volume 26 "-5.95%"

# Row 27 - PEPE
Set-Volume 27 "-8.11%"

# Row 28 - Binance-PegBSC-USD
Set-Volume 28 "+0.22%"

# Row 29 - RenderToken
Set-PriceText 29 "7.11"
Set-Volume 29 "-3.62%"

# Row 30 - Fetch.AI
Set-Volume 30 "-4.46%"

# Row 31 - InternetComputer(DFINITY)
Set-Volume 31 "-4.58%"

# Row 32 - PancakeSwap
Set-Volume 32 "-6.87%"

# Row 33 - USDe
Set-Volume 33 "-0.07%"

# Row 34 - Kaspa
Set-Volume 34 "-5.75%"

# Row 35 - RenzoRestakedETH
Set-Price 35 "3.340.26"
Set-Volume 35 "-5.22%"

# Row 36 - EthereumClassic
Set-PriceText 36 "22.31"
Set-Volume 36 "-4.24%"

# Row 37 - NEARProtocol
Set-PriceText 37 "5.33"
Set-Volume 37 "+0.03%"

# Row 38 & 39 - Monero and Aptos swap positions (Aptos now ranked above Monero)
$ws.Cells.Item(38, 2).Value = "Aptos"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-PriceText 38 "6.64"
Set-Volume 38 "-4.01%"

$ws.Cells.Item(39, 2).Value = "Monero"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-PriceText 39 "162.47"
Set-Volume 39 "-1.89%"

# Row 40 - ImmutableX
Set-PriceText 40 "1.47"
Set-Volume 40 "-4.60%"

# Row 41 - Hedera
Set-PriceText 41 "0.0748"
Set-Volume 41 "-4.82%"

# Row 42 - FirstDigitalUSD
Set-Volume 42 "+0.00%"

# Row 43 - OKB
Set-PriceText 43 "40.99"
Set-Volume 43 "-2.13%"

# Row 44 - Mantle
Set-Volume 44 "-8.04%"

# Row 45 - Filecoin
Set-PriceText 45 "4.20"
Set-Volume 45 "-5.03%"

# Row 46 - ONDO
Set-PriceText 46 "1.11"
Set-Volume 46 "-5.71%"

# Row 47 - Stacks
Set-Volume 47 "-6.75%"

# Row 48 - EnergySwap
Set-PriceText 48 "22.28"
Set-Volume 48 "-8.82%"

# Row 49 - Cosmos
Set-PriceText 49 "6.62"
Set-Volume 49 "-3.07%"

# Row 50 - SuiNetwork
Set-PriceText 50 "0.846"
Set-Volume 50 "-7.71%"

# Row 51 - Maker
Set-Price 51 "2.189.88"
Set-Volume 51 "-9.83%"
